$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# A leading apostrophe forces Excel to treat the assigned value as text (avoiding
# automatic numeric/date conversion for values such as '59.025.38' or '540.10'),
# and resetting the Style back to "Normal" keeps the cell formatting identical to
# the original (unstyled) inline-string cells.

$ws.Range("D2").Value = "'59.025.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -6.00%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.447.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -8.58%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'540.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.26%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'147.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -6.87%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.19%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.30%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.462.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -8.09%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.0991"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -6.41%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -2.10%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.30%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -4.53%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'2.884.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -8.39%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'23.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -9.55%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'58.899.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -6.05%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -6.33%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.516.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -6.02%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  -6.49%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -5.86%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'324.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.78%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.964"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.52%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'5.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -9.23%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'60.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.95%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.451"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -11.19%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -4.98%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -3.07%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -6.19%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -5.59%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.0₃0770"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -9.83%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'6.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.87%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -11.97%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -0.10%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'156.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.96%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -7.07%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'18.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -5.41%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -9.24%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D39").Value = "'316.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -10.19%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'5.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.69%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.838"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -11.75%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'36.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.48%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -6.95%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.28%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  -2.67%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  -3.04%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -5.92%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0525"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -6.30%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -5.29%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'121.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.50%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  -9.97%  "
$ws.Range("E51").Style = "Normal"
